$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newDate = Get-Date -Year 2023 -Month 10 -Day 4 -Hour 0 -Minute 0 -Second 0

for ($row = 2; $row -le 311; $row++) {
    $ws.Cells.Item($row, 3).Value = $newDate
}
